# Scheduled-runner refresh of market/profit figures in Sheets/Asura_Profits.xlsx
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) columns (H:N)
# for the leves whose backing Universalis market data changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 98 (Leve Item ID 36237)
if ($ws.Range("G98").Value2 -ne 36237) {
    Write-Output "skip ALC!98: Leve Item ID mismatch"
} else {
    $ws.Range("H98").Value = 5665.6177
    $ws.Range("I98").Value = 3754.3333
    $ws.Range("J98").Value = 20000.25
    $ws.Range("K98").Value = 3754.3333
    $ws.Range("L98").Value = 20000.25
    $ws.Range("M98").Value = -2256.3333
    $ws.Range("N98").Value = -22996.25
}

# ALC row 106 (Leve Item ID 19903)
if ($ws.Range("G106").Value2 -ne 19903) {
    Write-Output "skip ALC!106: Leve Item ID mismatch"
} else {
    $ws.Range("H106").Value = 1276.6666
    $ws.Range("I106").Value = 965
    $ws.Range("J106").Value = 1900
    $ws.Range("K106").Value = 965
    $ws.Range("L106").Value = 1900
    $ws.Range("M106").Value = -334
    $ws.Range("N106").Value = -3162
}

# ALC row 107 (Leve Item ID 27766)
if ($ws.Range("G107").Value2 -ne 27766) {
    Write-Output "skip ALC!107: Leve Item ID mismatch"
} else {
    $ws.Range("H107").Value = 100777.5
    $ws.Range("I107").Value = 143568.28
    $ws.Range("J107").Value = 932.3333
    $ws.Range("K107").Value = 143568.28
    $ws.Range("L107").Value = 932.3333
    $ws.Range("M107").Value = -141648.28
    $ws.Range("N107").Value = -4772.3333
}

# ALC row 112 (Leve Item ID 27960)
if ($ws.Range("G112").Value2 -ne 27960) {
    Write-Output "skip ALC!112: Leve Item ID mismatch"
} else {
    $ws.Range("H112").Value = 6708.077
    $ws.Range("J112").Value = 7470
    $ws.Range("L112").Value = 22410
    $ws.Range("N112").Value = -24626
}

# ALC row 122 (Leve Item ID 36237)
if ($ws.Range("G122").Value2 -ne 36237) {
    Write-Output "skip ALC!122: Leve Item ID mismatch"
} else {
    $ws.Range("H122").Value = 5665.6177
    $ws.Range("I122").Value = 3754.3333
    $ws.Range("J122").Value = 20000.25
    $ws.Range("K122").Value = 11262.9999
    $ws.Range("L122").Value = 60000.75
    $ws.Range("M122").Value = -8812.999899999999
    $ws.Range("N122").Value = -64900.75
}

# ALC row 129 (Leve Item ID 36115)
if ($ws.Range("G129").Value2 -ne 36115) {
    Write-Output "skip ALC!129: Leve Item ID mismatch"
} else {
    $ws.Range("H129").Value = 1374.95
    $ws.Range("J129").Value = 1427.7222
    $ws.Range("L129").Value = 4283.1666
    $ws.Range("N129").Value = -14283.1666
}

# ALC row 138 (Leve Item ID 44169)
if ($ws.Range("G138").Value2 -ne 44169) {
    Write-Output "skip ALC!138: Leve Item ID mismatch"
} else {
    $ws.Range("H138").Value = 3340367
    $ws.Range("I138").Value = 11769204
    $ws.Range("J138").Value = 8036.07
    $ws.Range("K138").Value = 35307612
    $ws.Range("L138").Value = 24108.21
    $ws.Range("M138").Value = -35302472
    $ws.Range("N138").Value = -34388.21
}

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2 (Leve Item ID 27713)
if ($ws.Range("G2").Value2 -ne 27713) {
    Write-Output "skip ARM!2: Leve Item ID mismatch"
} else {
    $ws.Range("H2").Value = 53652.527
    $ws.Range("I2").Value = 945.2727
    $ws.Range("J2").Value = 126125
    $ws.Range("K2").Value = 945.2727
    $ws.Range("L2").Value = 126125
    $ws.Range("M2").Value = -832.2727
    $ws.Range("N2").Value = -126351
}

# ARM row 32 (Leve Item ID 44147)
if ($ws.Range("G32").Value2 -ne 44147) {
    Write-Output "skip ARM!32: Leve Item ID mismatch"
} else {
    $ws.Range("H32").Value = 59102.61
    $ws.Range("I32").Value = 46386.684
    $ws.Range("J32").Value = 119503.25
    $ws.Range("K32").Value = 46386.684
    $ws.Range("L32").Value = 119503.25
    $ws.Range("M32").Value = -46099.684
    $ws.Range("N32").Value = -120077.25
}

# ARM row 45 (Leve Item ID 27714)
if ($ws.Range("G45").Value2 -ne 27714) {
    Write-Output "skip ARM!45: Leve Item ID mismatch"
} else {
    $ws.Range("H45").Value = 8014.5
    $ws.Range("I45").Value = 625.375
    $ws.Range("J45").Value = 17866.666
    $ws.Range("K45").Value = 625.375
    $ws.Range("L45").Value = 17866.666
    $ws.Range("M45").Value = -248.375
    $ws.Range("N45").Value = -18620.666
}

# ARM row 110 (Leve Item ID 27708)
if ($ws.Range("G110").Value2 -ne 27708) {
    Write-Output "skip ARM!110: Leve Item ID mismatch"
} else {
    $ws.Range("H110").Value = 1336.9166
    $ws.Range("I110").Value = 1264.3
    $ws.Range("J110").Value = 1700
    $ws.Range("K110").Value = 1264.3
    $ws.Range("L110").Value = 1700
    $ws.Range("M110").Value = 780.7
    $ws.Range("N110").Value = -5790
}

# ARM row 115 (Leve Item ID 27104)
if ($ws.Range("G115").Value2 -ne 27104) {
    Write-Output "skip ARM!115: Leve Item ID mismatch"
} else {
    $ws.Range("H115").Value = 20000
    $ws.Range("J115").Value = 20000
    $ws.Range("L115").Value = 20000
    $ws.Range("N115").Value = -23134
}

# ARM row 116 (Leve Item ID 27713)
if ($ws.Range("G116").Value2 -ne 27713) {
    Write-Output "skip ARM!116: Leve Item ID mismatch"
} else {
    $ws.Range("H116").Value = 53652.527
    $ws.Range("I116").Value = 945.2727
    $ws.Range("J116").Value = 126125
    $ws.Range("K116").Value = 945.2727
    $ws.Range("L116").Value = 126125
    $ws.Range("M116").Value = 1348.7273
    $ws.Range("N116").Value = -130713
}

# ARM row 122 (Leve Item ID 36168)
if ($ws.Range("G122").Value2 -ne 36168) {
    Write-Output "skip ARM!122: Leve Item ID mismatch"
} else {
    $ws.Range("H122").Value = 1476
    $ws.Range("I122").Value = 1296.6666
    $ws.Range("J122").Value = 2014
    $ws.Range("K122").Value = 3889.9998
    $ws.Range("L122").Value = 6042
    $ws.Range("M122").Value = -1439.9998
    $ws.Range("N122").Value = -10942
}

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3 (Leve Item ID 27713)
if ($ws.Range("G3").Value2 -ne 27713) {
    Write-Output "skip BSM!3: Leve Item ID mismatch"
} else {
    $ws.Range("H3").Value = 53652.527
    $ws.Range("I3").Value = 945.2727
    $ws.Range("J3").Value = 126125
    $ws.Range("K3").Value = 945.2727
    $ws.Range("L3").Value = 126125
    $ws.Range("M3").Value = -831.2727
    $ws.Range("N3").Value = -126353
}

# BSM row 107 (Leve Item ID 27706)
if ($ws.Range("G107").Value2 -ne 27706) {
    Write-Output "skip BSM!107: Leve Item ID mismatch"
} else {
    $ws.Range("H107").Value = 18087.281
    $ws.Range("I107").Value = 25104.363
    $ws.Range("J107").Value = 2649.7
    $ws.Range("K107").Value = 25104.363
    $ws.Range("L107").Value = 2649.7
    $ws.Range("M107").Value = -23184.363
    $ws.Range("N107").Value = -6489.7
}

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31 (Leve Item ID 44023)
if ($ws.Range("G31").Value2 -ne 44023) {
    Write-Output "skip CRP!31: Leve Item ID mismatch"
} else {
    $ws.Range("H31").Value = 5196
    $ws.Range("I31").Value = 4342.1665
    $ws.Range("J31").Value = 6476.75
    $ws.Range("K31").Value = 4342.1665
    $ws.Range("L31").Value = 6476.75
    $ws.Range("M31").Value = -4047.1665
    $ws.Range("N31").Value = -7066.75
}

# CRP row 34 (Leve Item ID 44023)
if ($ws.Range("G34").Value2 -ne 44023) {
    Write-Output "skip CRP!34: Leve Item ID mismatch"
} else {
    $ws.Range("H34").Value = 5196
    $ws.Range("I34").Value = 4342.1665
    $ws.Range("J34").Value = 6476.75
    $ws.Range("K34").Value = 4342.1665
    $ws.Range("L34").Value = 6476.75
    $ws.Range("M34").Value = -4140.1665
    $ws.Range("N34").Value = -6880.75
}

# CRP row 122 (Leve Item ID 36196)
if ($ws.Range("G122").Value2 -ne 36196) {
    Write-Output "skip CRP!122: Leve Item ID mismatch"
} else {
    $ws.Range("H122").Value = 2138
    $ws.Range("J122").Value = 2014
    $ws.Range("L122").Value = 6042
    $ws.Range("N122").Value = -10942
}

$ws = $wb.Worksheets.Item("CUL")
# CUL row 17 (Leve Item ID 4640)
if ($ws.Range("G17").Value2 -ne 4640) {
    Write-Output "skip CUL!17: Leve Item ID mismatch"
} else {
    $ws.Range("H17").Value = 946.7059
    $ws.Range("J17").Value = 1063.9286
    $ws.Range("L17").Value = 3191.7858
    $ws.Range("N17").Value = -3529.7858
}

# CUL row 34 (Leve Item ID 4749)
if ($ws.Range("G34").Value2 -ne 4749) {
    Write-Output "skip CUL!34: Leve Item ID mismatch"
} else {
    $ws.Range("H34").Value = 1012.1429
    $ws.Range("J34").Value = 1143.375
    $ws.Range("L34").Value = 3430.125
    $ws.Range("N34").Value = -3598.125
}

# CUL row 39 (Leve Item ID 4712)
if ($ws.Range("G39").Value2 -ne 4712) {
    Write-Output "skip CUL!39: Leve Item ID mismatch"
} else {
    $ws.Range("H39").Value = 2325.5715
    $ws.Range("J39").Value = 2975.8
    $ws.Range("L39").Value = 8927.400000000001
    $ws.Range("N39").Value = -9515.400000000001
}

# CUL row 55 (Leve Item ID 4733)
if ($ws.Range("G55").Value2 -ne 4733) {
    Write-Output "skip CUL!55: Leve Item ID mismatch"
} else {
    $ws.Range("H55").Value = 3562.3635
    $ws.Range("J55").Value = 3562.3635
    $ws.Range("L55").Value = 10687.0905
    $ws.Range("N55").Value = -11041.0905
}

# CUL row 131 (Leve Item ID 36060)
if ($ws.Range("G131").Value2 -ne 36060) {
    Write-Output "skip CUL!131: Leve Item ID mismatch"
} else {
    $ws.Range("H131").Value = 29416508
    $ws.Range("J131").Value = 32259992
    $ws.Range("L131").Value = 96779976
    $ws.Range("N131").Value = -96790056
}

$ws = $wb.Worksheets.Item("GSM")
# GSM row 102 (Leve Item ID 36169)
if ($ws.Range("G102").Value2 -ne 36169) {
    Write-Output "skip GSM!102: Leve Item ID mismatch"
} else {
    $ws.Range("H102").Value = 1751.0476
    $ws.Range("I102").Value = 1722.5714
    $ws.Range("J102").Value = 1808
    $ws.Range("K102").Value = 1722.5714
    $ws.Range("L102").Value = 1808
    $ws.Range("M102").Value = -100.5714
    $ws.Range("N102").Value = -5052
}

# GSM row 113 (Leve Item ID 27710)
if ($ws.Range("G113").Value2 -ne 27710) {
    Write-Output "skip GSM!113: Leve Item ID mismatch"
} else {
    $ws.Range("H113").Value = 1359.2727
    $ws.Range("I113").Value = 1094
    $ws.Range("J113").Value = 2066.6667
    $ws.Range("K113").Value = 1094
    $ws.Range("L113").Value = 2066.6667
    $ws.Range("M113").Value = 1076
    $ws.Range("N113").Value = -6406.6667
}

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22 (Leve Item ID 5277)
if ($ws.Range("G22").Value2 -ne 5277) {
    Write-Output "skip LTW!22: Leve Item ID mismatch"
} else {
    $ws.Range("H22").Value = 2908.5
    $ws.Range("I22").Value = 349
    $ws.Range("J22").Value = 3420.4
    $ws.Range("K22").Value = 349
    $ws.Range("L22").Value = 3420.4
    $ws.Range("M22").Value = -54
    $ws.Range("N22").Value = -4010.4
}

# LTW row 27 (Leve Item ID 5277)
if ($ws.Range("G27").Value2 -ne 5277) {
    Write-Output "skip LTW!27: Leve Item ID mismatch"
} else {
    $ws.Range("H27").Value = 2908.5
    $ws.Range("I27").Value = 349
    $ws.Range("J27").Value = 3420.4
    $ws.Range("K27").Value = 349
    $ws.Range("L27").Value = 3420.4
    $ws.Range("M27").Value = -242
    $ws.Range("N27").Value = -3634.4
}

# LTW row 42 (Leve Item ID 4333)
if ($ws.Range("G42").Value2 -ne 4333) {
    Write-Output "skip LTW!42: Leve Item ID mismatch"
} else {
    $ws.Range("H42").Value = 35000
    $ws.Range("J42").Value = 20000
    $ws.Range("L42").Value = 20000
    $ws.Range("N42").Value = -21126
}

# LTW row 46 (Leve Item ID 5282)
if ($ws.Range("G46").Value2 -ne 5282) {
    Write-Output "skip LTW!46: Leve Item ID mismatch"
} else {
    $ws.Range("H46").Value = 1264.3889
    $ws.Range("I46").Value = 1044.6666
    $ws.Range("J46").Value = 1374.25
    $ws.Range("K46").Value = 1044.6666
    $ws.Range("L46").Value = 1374.25
    $ws.Range("M46").Value = -856.6666
    $ws.Range("N46").Value = -1750.25
}

# LTW row 49 (Leve Item ID 4333)
if ($ws.Range("G49").Value2 -ne 4333) {
    Write-Output "skip LTW!49: Leve Item ID mismatch"
} else {
    $ws.Range("H49").Value = 35000
    $ws.Range("J49").Value = 20000
    $ws.Range("L49").Value = 20000
    $ws.Range("N49").Value = -20294
}

# LTW row 55 (Leve Item ID 5284)
if ($ws.Range("G55").Value2 -ne 5284) {
    Write-Output "skip LTW!55: Leve Item ID mismatch"
} else {
    $ws.Range("H55").Value = 225.09091
    $ws.Range("I55").Value = 273.25
    $ws.Range("J55").Value = 96.666664
    $ws.Range("K55").Value = 273.25
    $ws.Range("L55").Value = 96.666664
    $ws.Range("M55").Value = -100.25
    $ws.Range("N55").Value = -442.666664
}

# LTW row 58 (Leve Item ID 1728)
if ($ws.Range("G58").Value2 -ne 1728) {
    Write-Output "skip LTW!58: Leve Item ID mismatch"
} else {
    $ws.Range("H58").Value = 0
    $ws.Range("J58").Value = 0
    $ws.Range("L58").Value = 0
    $ws.Range("N58").Value = ""
}

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122 (Leve Item ID 36208)
if ($ws.Range("G122").Value2 -ne 36208) {
    Write-Output "skip WVR!122: Leve Item ID mismatch"
} else {
    $ws.Range("H122").Value = 22728888
    $ws.Range("I122").Value = 62501100
    $ws.Range("J122").Value = 1910
    $ws.Range("K122").Value = 187503300
    $ws.Range("L122").Value = 5730
    $ws.Range("M122").Value = -187500850
    $ws.Range("N122").Value = -10630
}
